$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "time_taken" in F1, copying the style from E1 (bold/centered/bordered header style)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill F2:F117 with the time_taken values for each data row
$timeTaken = @(
    "2021-10-05 10:50:21.025612",
    "2021-10-05 10:50:21.025623",
    "2021-10-05 10:50:21.025627",
    "2021-10-05 10:50:21.025630",
    "2021-10-05 10:50:21.025633",
    "2021-10-05 10:50:21.025635",
    "2021-10-05 10:50:21.025638",
    "2021-10-05 10:50:21.025641",
    "2021-10-05 10:50:21.025644",
    "2021-10-05 10:50:21.025646",
    "2021-10-05 10:50:21.025649",
    "2021-10-05 10:50:21.025651",
    "2021-10-05 10:50:21.025654",
    "2021-10-05 10:50:21.025656",
    "2021-10-05 10:50:21.025659",
    "2021-10-05 10:50:21.025661",
    "2021-10-05 10:50:21.025664",
    "2021-10-05 10:50:21.025667",
    "2021-10-05 10:50:21.025670",
    "2021-10-05 10:50:21.025672",
    "2021-10-05 10:50:21.025675",
    "2021-10-05 10:50:21.025677",
    "2021-10-05 10:50:21.025680",
    "2021-10-05 10:50:21.025682",
    "2021-10-05 10:50:21.025685",
    "2021-10-05 10:50:21.025688",
    "2021-10-05 10:50:21.025690",
    "2021-10-05 10:50:21.025693",
    "2021-10-05 10:50:21.025696",
    "2021-10-05 10:50:21.025698",
    "2021-10-05 10:50:21.025701",
    "2021-10-05 10:50:21.025703",
    "2021-10-05 10:50:21.025706",
    "2021-10-05 10:50:21.025709",
    "2021-10-05 10:50:21.025711",
    "2021-10-05 10:50:21.025714",
    "2021-10-05 10:50:21.025716",
    "2021-10-05 10:50:21.025719",
    "2021-10-05 10:50:21.025721",
    "2021-10-05 10:50:21.025724",
    "2021-10-05 10:50:21.025727",
    "2021-10-05 10:50:21.025729",
    "2021-10-05 10:50:21.025732",
    "2021-10-05 10:50:21.025735",
    "2021-10-05 10:50:21.025737",
    "2021-10-05 10:50:21.025740",
    "2021-10-05 10:50:21.025742",
    "2021-10-05 10:50:21.025745",
    "2021-10-05 10:50:21.025747",
    "2021-10-05 10:50:21.025749",
    "2021-10-05 10:50:21.025752",
    "2021-10-05 10:50:21.025754",
    "2021-10-05 10:50:21.025757",
    "2021-10-05 10:50:21.025760",
    "2021-10-05 10:50:21.025763",
    "2021-10-05 10:50:21.025765",
    "2021-10-05 10:50:21.025767",
    "2021-10-05 10:50:21.025770",
    "2021-10-05 10:50:21.025772",
    "2021-10-05 10:50:21.025775",
    "2021-10-05 10:50:21.025777",
    "2021-10-05 10:50:21.025780",
    "2021-10-05 10:50:21.025782",
    "2021-10-05 10:50:21.025785",
    "2021-10-05 10:50:21.025788",
    "2021-10-05 10:50:21.025791",
    "2021-10-05 10:50:21.025794",
    "2021-10-05 10:50:21.025796",
    "2021-10-05 10:50:21.025799",
    "2021-10-05 10:50:21.025802",
    "2021-10-05 10:50:21.025804",
    "2021-10-05 10:50:21.025807",
    "2021-10-05 10:50:21.025809",
    "2021-10-05 10:50:21.025812",
    "2021-10-05 10:50:21.025814",
    "2021-10-05 10:50:21.025817",
    "2021-10-05 10:50:21.025821",
    "2021-10-05 10:50:21.025824",
    "2021-10-05 10:50:21.025827",
    "2021-10-05 10:50:21.025829",
    "2021-10-05 10:50:21.025832",
    "2021-10-05 10:50:21.025834",
    "2021-10-05 10:50:21.025837",
    "2021-10-05 10:50:21.025839",
    "2021-10-05 10:50:21.025842",
    "2021-10-05 10:50:21.025844",
    "2021-10-05 10:50:21.025847",
    "2021-10-05 10:50:21.025849",
    "2021-10-05 10:50:21.025852",
    "2021-10-05 10:50:21.025854",
    "2021-10-05 10:50:21.025857",
    "2021-10-05 10:50:21.025859",
    "2021-10-05 10:50:21.025863",
    "2021-10-05 10:50:21.025866",
    "2021-10-05 10:50:21.025868",
    "2021-10-05 10:50:21.025871",
    "2021-10-05 10:50:21.025873",
    "2021-10-05 10:50:21.025876",
    "2021-10-05 10:50:21.025878",
    "2021-10-05 10:50:21.025881",
    "2021-10-05 10:50:21.025883",
    "2021-10-05 10:50:21.025886",
    "2021-10-05 10:50:21.025888",
    "2021-10-05 10:50:21.025891",
    "2021-10-05 10:50:21.025893",
    "2021-10-05 10:50:21.025896",
    "2021-10-05 10:50:21.025898",
    "2021-10-05 10:50:21.025901",
    "2021-10-05 10:50:21.025906",
    "2021-10-05 10:50:21.025909",
    "2021-10-05 10:50:21.025911",
    "2021-10-05 10:50:21.025914",
    "2021-10-05 10:50:21.025916",
    "2021-10-05 10:50:21.025919",
    "2021-10-05 10:50:21.025921",
    "2021-10-05 10:50:21.025924"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
